# Fixed shooting and bullets
# Fill in Week 40 (rows 131-136) and Week 41 (rows 139-144) attendance hours
# for Rief/David/Walter/Paige/Rick/Benjamin/Hadewij (columns B, D, E, G, H, I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 40 : Maandag..Vrijdag (rows 131-135) ---
# Maandag
$ws.Range("B131").Value = 4
$ws.Range("D131").Value = 4
$ws.Range("E131").Value = 4
$ws.Range("G131").Value = 4
$ws.Range("H131").Value = 4
$ws.Range("I131").Value = 4

# Dinsdag
$ws.Range("B132").Value = 6
$ws.Range("D132").Value = 6
$ws.Range("E132").Value = 6
$ws.Range("G132").Value = 6
$ws.Range("H132").Value = 6
$ws.Range("I132").Value = 6

# Woensdag
$ws.Range("B133").Value = 2
$ws.Range("D133").Value = 2
$ws.Range("E133").Value = 2
$ws.Range("G133").Value = 2
$ws.Range("H133").Value = 2
$ws.Range("I133").Value = 2

# Donderdag
$ws.Range("B134").Value = 4
$ws.Range("D134").Value = 4
$ws.Range("E134").Value = 4
$ws.Range("G134").Value = 4
$ws.Range("H134").Value = 4
$ws.Range("I134").Value = 4

# Vrijdag
$ws.Range("B135").Value = 4
$ws.Range("D135").Value = 4
$ws.Range("E135").Value = 4
$ws.Range("G135").Value = 4
$ws.Range("H135").Value = 4
$ws.Range("I135").Value = 4

# Totaal Game-Lab uren p/w (manually entered, not a formula)
$ws.Range("B136").Value = 20

# --- Week 41 : Maandag..Vrijdag (rows 139-143) ---
# Maandag
$ws.Range("B139").Value = 4
$ws.Range("D139").Value = 4
$ws.Range("E139").Value = 4
$ws.Range("G139").Value = 4
$ws.Range("H139").Value = 4
$ws.Range("I139").Value = 4

# Dinsdag
$ws.Range("B140").Value = 6
$ws.Range("D140").Value = 6
$ws.Range("E140").Value = 6
$ws.Range("G140").Value = 6
$ws.Range("H140").Value = 6
$ws.Range("I140").Value = 6

# Woensdag
$ws.Range("B141").Value = 2

# Donderdag
$ws.Range("B142").Value = 4

# Vrijdag
$ws.Range("B143").Value = 4

# Put the final selection where the author left it
$ws.Range("G140").Select()
